# Generate Report for Handoff
# For every file row whose status is "Handback transform failed" or
# "Ready for handoff", refresh the "Latest Handoff Date(time)" value to
# the timestamp of this report generation run.

$wb = $excel.ActiveWorkbook

$newOverviewDate  = "2016-22-17 20:22:33"
$newZhHandoffDate = "2016-03-17 20:22:30"
$newDeHandoffDate = "2016-03-17 20:22:33"

$statusesToRefresh = @("Handback transform failed", "Ready for handoff")

# --- Overview sheet: column B = status, column D = Latest Handoff Date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$lastRow = $wsOverview.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $status = $wsOverview.Cells.Item($r, 2).Text
    if ($statusesToRefresh -contains $status) {
        $wsOverview.Cells.Item($r, 4).Value = $newOverviewDate
    }
}

# --- zh-cn sheet: column C = status, column E = Latest Handoff Datetime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$lastRow = $wsZh.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $status = $wsZh.Cells.Item($r, 3).Text
    if ($statusesToRefresh -contains $status) {
        $wsZh.Cells.Item($r, 5).Value = $newZhHandoffDate
    }
}

# --- de-de sheet: column C = status, column E = Latest Handoff Datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
$lastRow = $wsDe.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $status = $wsDe.Cells.Item($r, 3).Text
    if ($statusesToRefresh -contains $status) {
        $wsDe.Cells.Item($r, 5).Value = $newDeHandoffDate
    }
}
